$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 174, shifting existing rows 174-275 down to 175-276
$ws.Rows("174").Insert()

# Populate the new row 174 with its values
$ws.Cells.Item(174, 1).Value = 10
$ws.Cells.Item(174, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(174, 3).Value = "La Araucanía"
$ws.Cells.Item(174, 4).Value = 44488
$ws.Cells.Item(174, 5).Value = 9
$ws.Cells.Item(174, 6).Value = 100112023
$ws.Cells.Item(174, 7).Value = "Brócoli"
$ws.Cells.Item(174, 8).Value = "Sin especificar"
$ws.Cells.Item(174, 9).Value = "Primera"
$ws.Cells.Item(174, 10).Value = 2500
$ws.Cells.Item(174, 11).Value = 800
$ws.Cells.Item(174, 12).Value = 800
$ws.Cells.Item(174, 13).Value = 800
$ws.Cells.Item(174, 14).Value = "$/unidad"
$ws.Cells.Item(174, 15).Value = "Región Metropolitana"
$ws.Cells.Item(174, 16).Value = 800
$ws.Cells.Item(174, 17).Value = 1
$ws.Cells.Item(174, 18).Value = "Hortaliza"

# Apply the date style (s="2") to the new D174 cell, matching the column's date formatting
$ws.Cells.Item(174, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
